# Aula 08 - "Contextualização" -> "Equações de Recorrência"
# (title slide + section-header slide)

$p = $ppt.ActivePresentation

# --- Slide 2: "Aula 08 <tab> Contextualização" title shape ---
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(3)
$tr2 = $sh2.TextFrame.TextRange

# Layout of the title text range:
#   1       -> line break (from <a:br>)
#   2-8     -> "Aula 08"
#   9       -> line break (from <a:br>)
#   10      -> tab character
#   11-26   -> "Contextualização"  (16 chars)
$titleRun = $tr2.Characters(11, 16)
$titleRun.Text = "Equações de Recorrência"

# --- Slide 3: "Contextualização" title shape ---
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(1)
$tr3 = $sh3.TextFrame.TextRange

# Original run text is "Contextualização" (16 chars): "Contextu" + "alização"
$part1 = $tr3.Characters(1, 8)
$part1.Text = "Equações"

$part2 = $tr3.Characters(9, 8)
$part2.Text = " de Recorrência"

# Split " de Recorrência" into " de " and "Recorrência" as separate runs
$recorrencia = $tr3.Characters(13, 11)
$recorrencia.Font.Bold = $true
